$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SQL queries stored in column B (and C2) joined tables using the
# surrogate key "id" (e.g. std.id / prt.id). They were updated to join on
# the actual natural-key id columns instead (e.g. std.study_id / prt.participant_id).
function Update-Query([string]$addr) {
    $s = $ws.Range($addr).Text
    $s = $s.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $s = $s.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $s = $s.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $s = $s.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $s = $s.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $s = $s.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
    $ws.Range($addr).Value = $s
}

Update-Query "B2"
Update-Query "C2"
Update-Query "B3"
Update-Query "B4"
Update-Query "B5"
Update-Query "B6"
Update-Query "B7"

# B5 previously carried a stray duplicate font style (12pt Calibri defined
# twice). Re-apply the same "wrap text, 12pt" look used by the sibling
# cells B6/B7 so the redundant style entry collapses away.
$ws.Range("B5").Font.Size = 12
$ws.Range("B5").WrapText = $true

# Selection/scroll moved from B2 to C7.
$ws.Activate()
$ws.Range("C7").Select() | Out-Null

Write-Host "Updated C3DC join queries for TC04 phs000469 AnatoSite-AdreGlndNOS."
